$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 12; this shifts rows 12..52 down to 13..53.
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with data (copy of what was the old row 12 pattern,
# updated with the new week's values).
$ws.Cells.Item(12, 1).Value = 7
$ws.Cells.Item(12, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(12, 3).Value = "Ñuble"
$ws.Cells.Item(12, 4).Value = 44575
$ws.Cells.Item(12, 4).NumberFormat = $ws.Cells.Item(13, 4).NumberFormat
$ws.Cells.Item(12, 5).Value = 16
$ws.Cells.Item(12, 6).Value = 100112022
$ws.Cells.Item(12, 7).Value = "Arveja Verde"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 120
$ws.Cells.Item(12, 11).Value = 24000
$ws.Cells.Item(12, 12).Value = 25000
$ws.Cells.Item(12, 13).Value = 24500
$ws.Cells.Item(12, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(12, 16).Value = 980
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"
